# Generate Report for Handback
# Refresh the Xliff generate / handoff / handback timestamps that the
# report-generation step stamps when it re-runs.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the efb15e1e... row.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-28 08:44:29"

# zh-cn sheet: Correspond Handoff / Handback datetimes for the efb15e1e... row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-28 08:44:24"
$wsZhCn.Range("K4").Value = "2016-08-28 08:45:05"

# de-de sheet: Correspond Handback datetime for the efb15e1e... row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K4").Value = "2016-08-28 08:45:14"
